$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the refreshed price strings (e.g. "1.000") are valid numeric
# literals, which Excel would otherwise silently coerce to a Number (losing
# the trailing zeros / formatting that the source feed uses). Force those
# specific cells to Text for the duration of the write, then restore the
# default "Normal" style so no stray number-format is left behind.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D18', 'D19', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D32', 'D33', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed coin data (price + 1h volume change), including the
# row 21/22 coin swap (Dai <-> WrappedliquidstakedEther2.0).
$ws.Range('D2').Value = '30.140.36'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.911.39'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('D5').Value = '0.7414'
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('D6').Value = '245.75'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.3101'
$ws.Range('E8').Value = '  -2.57%  '
$ws.Range('D9').Value = '26.56'
$ws.Range('E9').Value = '  -5.64%  '
$ws.Range('D10').Value = '0.06994'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('D11').Value = '0.08073'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').Value = '0.7724'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = '1.931.36'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '5.365'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '92.25'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').Value = '30.140.75'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '6.017'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = '0.000007883'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').Value = '240.60'
$ws.Range('E20').Value = '  -5.00%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.162.48'
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '7.128'
$ws.Range('E24').Value = '  +6.23%  '
$ws.Range('D25').Value = '9.439'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('D26').Value = '167.49'
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('D27').Value = '19.00'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('E28').Value = '  -2.79%  '
$ws.Range('D29').Value = '2.062'
$ws.Range('E29').Value = '  -6.72%  '
$ws.Range('D30').Value = '1.559'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').Value = '4.351'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('D33').Value = '4.107'
$ws.Range('E33').Value = '  -0.89%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').Value = '0.05171'
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('D36').Value = '0.7526'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').Value = '2.734'
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('D38').Value = '0.01954'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('D39').Value = '2.805'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '6.358'
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('D41').Value = '0.4525'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').Value = '74.80'
$ws.Range('E42').Value = '  -4.61%  '
$ws.Range('D43').Value = '2.000'
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('D44').Value = '0.8414'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').Value = '7.817'
$ws.Range('E46').Value = '  +2.94%  '
$ws.Range('D47').Value = '101.96'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = '9.958'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').Value = '2.065.71'
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('E51').Value = '  -2.06%  '

# Drop the temporary Text format now that the literal strings are stored,
# restoring each cell to its original (default) style.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
